# Update the "Förändrad" (changed) date column for rows 2-8 from 45174 to 45175.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 8; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45174) {
        $cell.Value = 45175
    }
}
